$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-289). Update it from 2023-10-03 (45202) to 2023-10-04 (45203).
$ws.Range("C2:C289").Value = 45203
